# Long_PFAS_Albuwell_Sample_051325.xlsx — "Add files via upload" re-save.
#
# The underlying data edits are on Sheet1, rows 86-89 (samples labelled
# H1-H4 in column A): the "Sample" label in column B is switched from the
# placeholder ladder-dilution names ("4 D6" / "8 D6") to the real sample
# names, the Dilution column (E) is corrected from 20 to 80, and a couple
# of Replicate values (F) are normalized. Row 89 (H4) is blanked out
# entirely (it becomes an empty data row, matching the blank rows already
# present for H5-H12 further down).
#
# Because "4 D6" and "8 D6" end up unused anywhere in the sheet once these
# edits are made, Excel's shared-string table naturally shrinks by those
# two entries on save (count 146->144, uniqueCount 112->110), which is
# also reflected in the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 86 (Cell H1): rename sample, fix dilution ---
$ws.Range("B86").Value = "3920_Final"
$ws.Range("E86").Value = 80

# --- Row 87 (Cell H2): rename sample, fix dilution + replicate ---
$ws.Range("B87").Value = "3932_Final"
$ws.Range("E87").Value = 80
$ws.Range("F87").Value = 1

# --- Row 88 (Cell H3): rename sample, fix dilution ---
$ws.Range("B88").Value = "3936_Final"
$ws.Range("E88").Value = 80

# --- Row 89 (Cell H4): clear out the Sample/Type/Dilution/Replicate data ---
$ws.Range("B89:F89").ClearContents()

# --- Update the on-screen view to match where the author left the cursor ---
$win = $excel.ActiveWindow
$win.ScrollRow = 38
$win.ScrollColumn = 1
$ws.Range("C80").Select()

Write-Output "Applied Long_PFAS_Albuwell_Sample_051325 edits (rows 86-89 + view state)."
